$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking scrape refresh: update Price (D) and Volume(1h) (E) text cells.
# Both columns store plain text (e.g. "538.79", "  +0.88%  "), never real
# numbers/percentages, so for values that would otherwise auto-parse as a
# number we briefly force a Text format, then restore the default "Normal"
# style so the cell keeps the workbook's original (unstyled) look.

$ws.Range("D2").Value = "59.508.78"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.606.93"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "3.062.42"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "59.420.10"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "2.570.09"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.409"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("D28").Value = "0.0₃0747"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +6.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.847"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.828"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").Value = "1.952.19"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.46%  "
